$d = $word.ActiveDocument

# 1. Merge ": " into a single run (cosmetic - Find/Replace text is unaffected
#    visually, but ensures the text still reads correctly).
$d.Content.Find.Execute("General Appeal of the Course:  ", $true, $false, $false, $false, $false, $true, 1, $false, "General Appeal of the Course: ", 2) | Out-Null

# 2. Spelling corrections inside the long paragraph.
$d.Content.Find.Execute("Institure", $true, $false, $false, $false, $false, $true, 1, $false, "Institute", 2) | Out-Null
$d.Content.Find.Execute("certifcation", $true, $false, $false, $false, $false, $true, 1, $false, "certification", 2) | Out-Null
$d.Content.Find.Execute("reqire", $true, $false, $false, $false, $false, $true, 1, $false, "require", 2) | Out-Null
